$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Macroalgae / Bloom phenology" description cell (E10).
# Mirror the formatting already used by D10 (wrapped, vertically centered
# text) and nudge the fill so the cell records its own style entry,
# matching the workbook's "Description" column styling for this row.
$ws.Range("E10").Value = "-"
$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Interior.ColorIndex = -4142

# Move the active selection to B12, as left by the edit.
$ws.Range("B12").Select()
